$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 619.3421170982774
$ws.Range("D2").Value = 10847.19484611139
